$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws "D2" "60.462.71"
Set-TextValue $ws "E2" "  +5.76%  "
Set-TextValue $ws "D3" "3.275.21"
Set-TextValue $ws "E3" "  +0.78%  "
Set-TextValue $ws "E4" "  +0.17%  "
Set-TextValue $ws "D5" "406.55"
Set-TextValue $ws "E5" "  +2.79%  "
Set-TextValue $ws "D6" "110.57"
Set-TextValue $ws "E6" "  +2.40%  "
Set-TextValue $ws "D7" "3.268.76"
Set-TextValue $ws "E7" "  +0.73%  "
Set-TextValue $ws "D8" "0.562"
Set-TextValue $ws "E8" "  -4.32%  "
Set-TextValue $ws "D9" "0.999"
Set-TextValue $ws "D10" "0.614"
Set-TextValue $ws "E10" "  -1.77%  "
Set-TextValue $ws "D11" "0.112"
Set-TextValue $ws "E11" "  +12.98%  "
Set-TextValue $ws "D12" "38.31"
Set-TextValue $ws "E12" "  -2.15%  "
Set-TextValue $ws "E13" "  -0.33%  "
Set-TextValue $ws "D14" "3.811.64"
Set-TextValue $ws "E14" "  +1.35%  "
Set-TextValue $ws "D15" "8.06"
Set-TextValue $ws "E15" "  -2.15%  "
Set-TextValue $ws "D16" "18.82"
Set-TextValue $ws "E16" "  -1.50%  "
Set-TextValue $ws "D17" "3.335.21"
Set-TextValue $ws "E17" "  +2.40%  "
Set-TextValue $ws "D18" "60.500.69"
Set-TextValue $ws "E18" "  +6.16%  "
Set-TextValue $ws "D19" "0.975"
Set-TextValue $ws "E19" "  -5.48%  "
Set-TextValue $ws "E20" "  -2.96%  "
Set-TextValue $ws "D21" "0.0000113"
Set-TextValue $ws "E21" "  +0.06%  "
Set-TextValue $ws "E22" "  -2.84%  "
Set-TextValue $ws "D23" "12.40"
Set-TextValue $ws "E23" "  -4.19%  "
Set-TextValue $ws "D24" "294.98"
Set-TextValue $ws "E24" "  -0.42%  "
Set-TextValue $ws "D25" "72.76"
Set-TextValue $ws "E25" "  -2.04%  "
Set-TextValue $ws "D26" "3.06"
Set-TextValue $ws "E26" "  -3.62%  "
Set-TextValue $ws "D27" "28.86"
Set-TextValue $ws "E27" "  +3.18%  "
Set-TextValue $ws "E28" "  -2.18%  "
Set-TextValue $ws "E29" "  +2.42%  "
Set-TextValue $ws "D30" "7.28"
Set-TextValue $ws "E30" "  +0.44%  "
Set-TextValue $ws "D31" "7.39"
Set-TextValue $ws "E31" "  -2.63%  "
Set-TextValue $ws "B32" "Hedera"
Set-TextValue $ws "C32" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws "D32" "0.111"
Set-TextValue $ws "E32" "  +2.40%  "
Set-TextValue $ws "B33" "Dai"
Set-TextValue $ws "C33" "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws "D33" "1.00"
Set-TextValue $ws "E33" "  +0.02%  "
Set-TextValue $ws "D34" "11.05"
Set-TextValue $ws "E34" "  -2.38%  "
Set-TextValue $ws "D35" "2.43"
Set-TextValue $ws "E35" "  +14.15%  "
Set-TextValue $ws "D36" "38.53"
Set-TextValue $ws "E36" "  -1.72%  "
Set-TextValue $ws "D37" "0.0473"
Set-TextValue $ws "E37" "  -1.72%  "
Set-TextValue $ws "D38" "52.05"
Set-TextValue $ws "E38" "  +0.87%  "
Set-TextValue $ws "D39" "1.00"
Set-TextValue $ws "E39" "  +0.35%  "
Set-TextValue $ws "D40" "3.09"
Set-TextValue $ws "E40" "  +5.39%  "
Set-TextValue $ws "D41" "3.26"
Set-TextValue $ws "E41" "  -6.45%  "
Set-TextValue $ws "D42" "134.09"
Set-TextValue $ws "E42" "  -0.31%  "
Set-TextValue $ws "B43" "Stellar"
Set-TextValue $ws "C43" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws "D43" "0.119"
Set-TextValue $ws "E43" "  -3.45%  "
Set-TextValue $ws "B44" "TheGraph"
Set-TextValue $ws "C44" "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws "D44" "0.284"
Set-TextValue $ws "E44" "  +1.46%  "
Set-TextValue $ws "D45" "1.86"
Set-TextValue $ws "E45" "  -1.44%  "
Set-TextValue $ws "D46" "16.08"
Set-TextValue $ws "E46" "  -5.94%  "
Set-TextValue $ws "D47" "3.70"
Set-TextValue $ws "E47" "  -6.42%  "
Set-TextValue $ws "E48" "  +2.21%  "
Set-TextValue $ws "D49" "20.72"
Set-TextValue $ws "E49" "  -6.70%  "
Set-TextValue $ws "D50" "2.099.72"
Set-TextValue $ws "E50" "  -2.66%  "
Set-TextValue $ws "B51" "ApeXProtocol"
Set-TextValue $ws "C51" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws "D51" "2.34"
Set-TextValue $ws "E51" "  +0.10%  "
